$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: fill in new Week 3 task assignments
$ws.Range("B11").Value = "Implement AI and HUD"
$ws.Range("C11").Value = "Implement tutorial"
$ws.Range("D11").Value = "TCR's: Resolution on start, minimize, fullscreen, "
$ws.Range("E11").Value = "Finalize zilch/work on polishing any issues "

# Row 14: new row of "Work on your guide" entries
$ws.Range("B14").Value = "Work on your guide"
$ws.Range("C14").Value = "Work on your guide"
$ws.Range("D14").Value = "Work on your guide"
$ws.Range("E14").Value = "Work on your guide"

# Row 12: clear old leftover tasks, keep only D12 with new text
$ws.Range("B12").Value = $null
$ws.Range("C12").Value = $null
$ws.Range("D12").Value = "editor fixes, work on menus, "
$ws.Range("E12").Value = $null
$ws.Range("F12").Value = $null

# Row 13: clear (was "editor fixes, work on menus" in D13)
$ws.Range("D13").Value = $null

# Last new task text for row 11
$ws.Range("F11").Value = "Implement an auto play "

# Row 16: remove the "Week 4" label and "Audio Guide" entry
$ws.Range("A16").Value = $null
$ws.Range("F16").Value = $null

# Row 17: remove the "Self-Play" entry
$ws.Range("D17").Value = $null

# Update the selected cell to match the new active selection
$ws.Range("F14").Select()
